# Commit #27 "Encore le face tracker" — content edits on the "Solution"
# slides (sensors / Kinect section) and the "Solution (organisation)"
# slide (explanation of interactions).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3 — "Solution": Content Placeholder 2
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$tr3 = $slide3.Shapes.Item(2).TextFrame.TextRange

# "L'utilisateur à le plein contrôle" -> split into two runs with
# identical formatting ("L'utilisateur à le plein " + "contrôle").
$paraControle = $tr3.Paragraphs(10, 1)
$sub = $paraControle.Characters(26, 8)
$sub.Text = $sub.Text

# New paragraphs appended at the end of the placeholder:
#   Capteur                (lvl 0)
#     Kinect                (lvl 1)
#     Intel                 (lvl 1)
#     PS Eye (Peut-être)    (lvl 1)
$newP = $tr3.InsertAfter([char]13 + "Capteur")
$cnt = $tr3.Paragraphs().Count
$tr3.Paragraphs($cnt, 1).IndentLevel = 1

$newP = $tr3.InsertAfter([char]13 + "Kinect")
$cnt = $tr3.Paragraphs().Count
$tr3.Paragraphs($cnt, 1).IndentLevel = 2

$newP = $tr3.InsertAfter([char]13 + "Intel")
$cnt = $tr3.Paragraphs().Count
$tr3.Paragraphs($cnt, 1).IndentLevel = 2

$newP = $tr3.InsertAfter([char]13 + "PS Eye (Peut-" + [char]0xEA + "tre)")
$cnt = $tr3.Paragraphs().Count
$tr3.Paragraphs($cnt, 1).IndentLevel = 2

# ---------------------------------------------------------------------
# Slide 4 — "Solution (organisation)": Content Placeholder 2
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$tr4 = $slide4.Shapes.Item(2).TextFrame.TextRange

# "Explication du jeu et de son fonctionnement" -> split into two runs
# with identical formatting.
$paraFonct = $tr4.Paragraphs(2, 1)
$subF = $paraFonct.Characters(30, 14)
$subF.Text = $subF.Text

# New paragraph right after it: "Explication des interactions (François)"
# at the same indent level (lvl 1).
$paraFonct2 = $tr4.Paragraphs(2, 1)
$newP2 = $paraFonct2.InsertAfter([char]13 + "Explication des interactions (Fran" + [char]0xE7 + "ois)")
$cnt4 = $tr4.Paragraphs().Count
$tr4.Paragraphs(3, 1).IndentLevel = 2
